# Auto-generated edit script for adc_sbb_within_100.docx
$d = $word.ActiveDocument

# Update the date/title line at the top of the document
$d.Content.Find.Execute("2026-01-05 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-06 Tuesday", 2) | Out-Null

# Update each arithmetic expression cell in the 20x5 table, by (row, col) position
$t = $d.Tables(1)

$newValues = @(
    @("51-24=", "34-15=", "15+8=", "85-29=", "37+48="),
    @("8+26=", "32-16=", "90-11=", "66+7=", "35+29="),
    @("77+5=", "31-18=", "62-49=", "78-49=", "81-7="),
    @("62-45=", "86-79=", "76-17=", "48+7=", "49+9="),
    @("54-8=", "91-29=", "54-5=", "7+67=", "51-15="),
    @("81-29=", "83-79=", "81-14=", "92-8=", "93-39="),
    @("70-52=", "36-18=", "79+5=", "51-8=", "17+55="),
    @("5+79=", "63+8=", "60-44=", "91-22=", "36-19="),
    @("51-18=", "7+39=", "9+73=", "26+35=", "70-1="),
    @("37+56=", "30-4=", "95-87=", "85-27=", "61-14="),
    @("90-34=", "61-12=", "56+7=", "81-53=", "47+47="),
    @("75+17=", "56+15=", "12+49=", "59+29=", "23+19="),
    @("90-87=", "73-4=", "80-37=", "7+77=", "92-23="),
    @("47+19=", "62-49=", "78+7=", "62-58=", "66-19="),
    @("94-47=", "67-39=", "18+8=", "35+56=", "81-53="),
    @("7+36=", "40-35=", "20-15=", "68-19=", "46+38="),
    @("86-9=", "58+37=", "65-36=", "24+9=", "84-38="),
    @("77-68=", "74-55=", "9+5=", "77-9=", "73-55="),
    @("72-49=", "26+48=", "36+57=", "30-19=", "42-18="),
    @("81-64=", "72+19=", "44+7=", "18+24=", "82-67=")
)

for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$r - 1][$c - 1]
    }
}

Write-Host "Done updating title and" ($t.Rows.Count * $t.Columns.Count) "table cells."
